$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Age"

$ws.Range("A3").Value = "sub 3058"
$ws.Range("B3").Value = 61
$ws.Range("C3").Value = 1422.428571428571
$ws.Range("D3").Value = 7.383171065316551
$ws.Range("E3").Value = 7.434027401184465

$ws.Range("A4").Value = "sub 2917"
$ws.Range("B4").Value = 59
$ws.Range("C4").Value = 450.8571428571428
$ws.Range("D4").Value = 6.92581662006612
$ws.Range("E4").Value = 6.830068449663596

$ws.Range("A5").Value = "sub 3104"
$ws.Range("B5").Value = 50
$ws.Range("C5").Value = 405.5714285714286
$ws.Range("D5").Value = 6.878098400502704
$ws.Range("E5").Value = 6.876968965035582

$ws.Range("A6").Value = "sub 3889"
$ws.Range("B6").Value = 62
$ws.Range("C6").Value = 202.7142857142857
$ws.Range("D6").Value = 7.34555419197442
$ws.Range("E6").Value = 7.434867293404128

$ws.Range("A7").Value = "sub 3642"
$ws.Range("B7").Value = 60
$ws.Range("C7").Value = 275
$ws.Range("D7").Value = 7.006946805494435
$ws.Range("E7").Value = 7.120774063177426

$ws.Range("A8").Value = "sub 3035"
$ws.Range("B8").Value = 67
$ws.Range("C8").Value = 843.2857142857143
$ws.Range("D8").Value = 7.199048679873495
$ws.Range("E8").Value = 7.344855788232741

$ws.Range("A9").Value = "sub 4182"
$ws.Range("B9").Value = 53
$ws.Range("C9").Value = 146.2857142857143
$ws.Range("D9").Value = 7.246077654891511
$ws.Range("E9").Value = 7.213574505912566

$ws.Range("A10").Value = "sub 4137"
$ws.Range("B10").Value = 70
$ws.Range("C10").Value = 162.2857142857143
$ws.Range("D10").Value = 6.779599736754784
$ws.Range("E10").Value = 6.820765719936312

$ws.Range("A11").Value = "sub 3583"
$ws.Range("B11").Value = 51
$ws.Range("C11").Value = 481.5714285714286
$ws.Range("D11").Value = 7.04089640854805
$ws.Range("E11").Value = 7.031325458341337

$ws.Range("A12").Value = "sub 4281"
$ws.Range("B12").Value = 46
$ws.Range("C12").Value = 128.5714285714286
$ws.Range("D12").Value = 6.938525114908565
$ws.Range("E12").Value = 7.064423196070866

$ws.Range("A13").Value = "sub 4191"
$ws.Range("B13").Value = 55
$ws.Range("C13").Value = 131.4285714285714
$ws.Range("D13").Value = 6.528340210711074
$ws.Range("E13").Value = 6.579365218778081

$ws.Range("A14").Value = "sub 3201"
$ws.Range("B14").Value = 51
$ws.Range("C14").Value = 436.4285714285714
$ws.Range("D14").Value = 7.064064798997272
$ws.Range("E14").Value = 7.064543698287891

$ws.Range("A15").Value = "sub 4275"
$ws.Range("B15").Value = 55
$ws.Range("C15").Value = 121.5714285714286
$ws.Range("D15").Value = 6.622892052670723
$ws.Range("E15").Value = 6.720497666077411

$ws.Range("A16").Value = "sub 4170"
$ws.Range("B16").Value = 58
$ws.Range("C16").Value = 308
$ws.Range("D16").Value = 7.017351638359888
$ws.Range("E16").Value = 7.152626773849482

$ws.Range("A17").Value = "sub 3676"
$ws.Range("B17").Value = 63
$ws.Range("C17").Value = 280.8571428571428
$ws.Range("D17").Value = 7.036606214789604
$ws.Range("E17").Value = 7.054626929132155

$ws.Range("A18").Value = "sub 4431"
$ws.Range("B18").Value = 71
$ws.Range("C18").Value = 96.71428571428571
$ws.Range("D18").Value = 7.020818116522912
$ws.Range("E18").Value = 7.068947435849092

$ws.Range("A19").Value = "sub 4294"
$ws.Range("B19").Value = 67
$ws.Range("C19").Value = 107.5714285714286
$ws.Range("D19").Value = 7.214559600266331
$ws.Range("E19").Value = 7.352766992376172

$ws.Range("A20").Value = "sub 3650"
$ws.Range("B20").Value = 62
$ws.Range("C20").Value = 331.5714285714286
$ws.Range("D20").Value = 7.071359391813155
$ws.Range("E20").Value = 7.304801206660142

$ws.Range("A21").Value = "sub 4140"
$ws.Range("B21").Value = 44
$ws.Range("C21").Value = 161.7142857142857
$ws.Range("D21").Value = 6.975010277606088
$ws.Range("E21").Value = 7.074018353204599

$ws.Range("A22").Value = "sub 3912"
$ws.Range("B22").Value = 55
$ws.Range("C22").Value = 229.8571428571429
$ws.Range("D22").Value = 7.118260683603749
$ws.Range("E22").Value = 7.106964680000716

$ws.Range("A23").Value = "sub 4051"
$ws.Range("B23").Value = 59
$ws.Range("C23").Value = 196.7142857142857
$ws.Range("D23").Value = 7.058402059460396
$ws.Range("E23").Value = 7.265104995210833

$ws.Range("A24").Value = "sub 3286"
$ws.Range("B24").Value = 76
$ws.Range("C24").Value = 390.1428571428572
$ws.Range("D24").Value = 6.947203699001764
$ws.Range("E24").Value = 7.071932059675354

$ws.Range("A25").Value = "sub 4208"
$ws.Range("B25").Value = 51
$ws.Range("C25").Value = 137.2857142857143
$ws.Range("D25").Value = 6.844582035474438
$ws.Range("E25").Value = 6.870261074554556

$ws.Range("A26").Value = "sub 1536"
$ws.Range("B26").Value = 69
$ws.Range("C26").Value = 1061
$ws.Range("D26").Value = 6.884511820778181
$ws.Range("E26").Value = 6.907893730935685

$ws.Range("A27").Value = "sub 4439"
$ws.Range("B27").Value = 56
$ws.Range("C27").Value = 79.14285714285714
$ws.Range("D27").Value = 7.240973929106727
$ws.Range("E27").Value = 7.266618338198525

$ws.Range("A28").Value = "sub 4466"
$ws.Range("B28").Value = 65
$ws.Range("C28").Value = 95.85714285714286
$ws.Range("D28").Value = 6.800857123428436
$ws.Range("E28").Value = 6.759607939389943

$ws.Range("A29").Value = "sub 4504"
$ws.Range("B29").Value = 55
$ws.Range("C29").Value = 78.14285714285714
$ws.Range("D29").Value = 7.17336214378206
$ws.Range("E29").Value = 7.183199988711299

$ws.Range("A30").Value = "sub 3887"
$ws.Range("B30").Value = 68
$ws.Range("C30").Value = 310.4285714285714
$ws.Range("D30").Value = 6.836989700151324
$ws.Range("E30").Value = 6.922828603479615

$ws.Range("A31").Value = "sub 4391"
$ws.Range("B31").Value = 61
$ws.Range("C31").Value = 101
$ws.Range("D31").Value = 6.89020938905185
$ws.Range("E31").Value = 7.040639403350221

$ws.Range("A32").Value = "sub 4546"
$ws.Range("B32").Value = 40
$ws.Range("C32").Value = 123.5714285714286
$ws.Range("D32").Value = 6.714085246818761
$ws.Range("E32").Value = 6.679859763459568

$ws.Range("A33").Value = "sub 4227"
$ws.Range("B33").Value = 38
$ws.Range("C33").Value = 186.7142857142857
$ws.Range("D33").Value = 6.556386429835174
$ws.Range("E33").ClearContents()

$ws.Range("A34").Value = "sub 4532"
$ws.Range("B34").Value = 63
$ws.Range("C34").Value = 89.42857142857143
$ws.Range("D34").Value = 7.038688554197772
$ws.Range("E34").Value = 7.054795021631452

$ws.Range("A35").Value = "sub 3299"
$ws.Range("B35").Value = 60
$ws.Range("C35").Value = 509.2857142857143
$ws.Range("D35").Value = 7.003880553503734
$ws.Range("E35").Value = 7.121080051765402

$ws.Range("A36").Value = "sub 4720"
$ws.Range("B36").Value = 59
$ws.Range("C36").Value = 43.28571428571428
$ws.Range("D36").Value = 6.9603032546778
$ws.Range("E36").Value = 7.074549500824056

$ws.Range("A37").Value = "sub 4607"
$ws.Range("B37").Value = 63
$ws.Range("C37").Value = 60.42857142857143
$ws.Range("D37").Value = 7.233524545965555
$ws.Range("E37").Value = 7.326711135958234

$ws.Range("A38").Value = "sub 4728"
$ws.Range("B38").Value = 43
$ws.Range("C38").Value = 34.42857142857143
$ws.Range("D38").Value = 6.993156859826387
$ws.Range("E38").Value = 7.002732069885225

$ws.Range("A39").Value = "sub 4769"
$ws.Range("B39").Value = 54
$ws.Range("C39").Value = 31.71428571428572
$ws.Range("D39").Value = 7.38145000440992
$ws.Range("E39").Value = 7.457650505583721

$ws.Range("A40").Value = "sub 4777"
$ws.Range("B40").Value = 51
$ws.Range("C40").Value = 28.71428571428572
$ws.Range("D40").Value = 6.882661353482145
$ws.Range("E40").Value = 7.007937705564125

$ws.Range("A41").Value = "sub 3396"
$ws.Range("B41").Value = 54
$ws.Range("C41").Value = 337.2857142857143
$ws.Range("D41").Value = 6.886312736330942
$ws.Range("E41").Value = 6.954748958644287

$ws.Range("A42").Value = "sub 4467"
$ws.Range("B42").Value = 60
$ws.Range("C42").Value = 31.28571428571428
$ws.Range("D42").Value = 7.195062281103835
$ws.Range("E42").Value = 7.124970950253855

$ws.Range("A43").Value = "sub 3392"
$ws.Range("B43").Value = 53
$ws.Range("C43").Value = 331.1428571428572
$ws.Range("D43").Value = 7.048890045500932
$ws.Range("E43").Value = 7.091926217511165

$ws.Range("A44").Value = "sub 3960"
$ws.Range("B44").Value = 79
$ws.Range("C44").Value = 239.7142857142857
$ws.Range("D44").Value = 7.272919090316775
$ws.Range("E44").Value = 7.365944104475781

$ws.Range("A45").Value = "sub 3154"
$ws.Range("B45").Value = 30
$ws.Range("C45").Value = 443.1428571428572
$ws.Range("D45").Value = 7.096204594757248
$ws.Range("E45").Value = 6.999666939315914

$ws.Range("A46").Value = "sub 3994"
$ws.Range("B46").Value = 55
$ws.Range("C46").Value = 179.8571428571429
$ws.Range("D46").Value = 7.118088235006313
$ws.Range("E46").Value = 7.115537975879147

$ws.Range("A47").Value = "sub 4086"
$ws.Range("B47").Value = 59
$ws.Range("C47").Value = 141.7142857142857
$ws.Range("D47").Value = 7.065382666671098
$ws.Range("E47").Value = 7.211965264109365

$ws.Range("A48").Value = "sub 4002"
$ws.Range("B48").Value = 69
$ws.Range("C48").Value = 185.2857142857143
$ws.Range("D48").Value = 7.295132844486576
$ws.Range("E48").Value = 7.307564638507666

$ws.Range("A49").Value = "sub 2522"
$ws.Range("B49").Value = 40
$ws.Range("C49").Value = 555.8571428571429
$ws.Range("D49").Value = 6.899297846806425
$ws.Range("E49").Value = 6.884580838538393

$ws.Range("A50").Value = "sub 3364"
$ws.Range("B50").Value = 34
$ws.Range("C50").Value = 334.4285714285714
$ws.Range("D50").Value = 7.091708295549481
$ws.Range("E50").Value = 7.019923080878011

$ws.Range("A51").Value = "sub 2884"
$ws.Range("B51").Value = 55
$ws.Range("C51").Value = 516.8571428571429
$ws.Range("D51").Value = 6.974732413069717
$ws.Range("E51").Value = 7.02803811151775

$ws.Range("A52").Value = "sub 2998"
$ws.Range("B52").Value = 59
$ws.Range("C52").Value = 469.2857142857143
$ws.Range("D52").Value = 7.316782896088275
$ws.Range("E52").Value = 7.537703125680994
